# This workbook holds one record (2 rows: "Primera"/"Segunda" quality)
# per reporting date, starting at row 311 and going through row 438.
# A brand-new weekly record needs to be inserted at the very top of that
# block (row 311), pushing every existing record down by one record
# (2 rows). That naturally creates the two new trailing rows (439/440)
# that used to be 437/438, and leaves the dimension as A1:R440.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 311), shifting
# all the existing records (311..438) down to (313..440).
$ws.Rows.Item(311).Insert()
$ws.Rows.Item(312).Insert()

# Fill in the new record in the freshly inserted rows 311 (Primera) and
# 312 (Segunda). Every column besides D/J/K/L/M/P is identical across all
# records in this sheet.
$ws.Range("A311").Value = 8
$ws.Range("B311").Value = "Terminal La Palmera de La Serena"
$ws.Range("C311").Value = "Coquimbo"
$ws.Range("D311").Value = 44755
$ws.Range("E311").Value = 4
$ws.Range("F311").Value = 100112017
$ws.Range("G311").Value = "Apio"
$ws.Range("H311").Value = "Americana (o)"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 2400
$ws.Range("K311").Value = 8500
$ws.Range("L311").Value = 9000
$ws.Range("M311").Value = 8750
$ws.Range("N311").Value = "`$/docena de matas"
$ws.Range("O311").Value = "Provincia del Elquí"
$ws.Range("P311").Value = 1458
$ws.Range("Q311").Value = 6
$ws.Range("R311").Value = "Hortaliza"

$ws.Range("A312").Value = 8
$ws.Range("B312").Value = "Terminal La Palmera de La Serena"
$ws.Range("C312").Value = "Coquimbo"
$ws.Range("D312").Value = 44755
$ws.Range("E312").Value = 4
$ws.Range("F312").Value = 100112017
$ws.Range("G312").Value = "Apio"
$ws.Range("H312").Value = "Americana (o)"
$ws.Range("I312").Value = "Segunda"
$ws.Range("J312").Value = 1400
$ws.Range("K312").Value = 7500
$ws.Range("L312").Value = 8000
$ws.Range("M312").Value = 7750
$ws.Range("N312").Value = "`$/docena de matas"
$ws.Range("O312").Value = "Provincia del Elquí"
$ws.Range("P312").Value = 1292
$ws.Range("Q312").Value = 6
$ws.Range("R312").Value = "Hortaliza"
